$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.582.10"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.638.53"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.63"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.09"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.173"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.637.97"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.94"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.115.94"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.326.76"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.82"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.625.87"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.10"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.87"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.58"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.45"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.767.51"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0952"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "494.05"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.52"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").Value = "  +6.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.14"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.90"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.326"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.02"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.87"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.545"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.63"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.609"
$ws.Range("E51").Value = "  +0.57%  "
